# Updated Table of Contents
# Adds a new row (23) to the Table1 listing for the
# "Pairs Violating BST property" problem.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry: Index 23, Problem Statement, TYPE 1, Difficulty
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = "Pairs Violating BST property"
$ws.Cells.Item(25, 3).Value = "Binary Search Tree"
$ws.Cells.Item(25, 5).Value = "Medium"

# Move the active selection the way Excel leaves it after typing the
# last entry and pressing Enter (drops to the row below, column E).
$ws.Range("E26").Select() | Out-Null
